# Update 29th June entries: add a new completed log-book row (Sno 31) for
# 2022-06-29 and refresh the selected cell, per commit "updated 29th june entries".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logBook")

# Row 32 is currently a blank gap row between the last entry (row 31) and
# the "Total Hours" row (row 33) - the sheet dimension already reserves it
# (A1:G33), so this is a plain fill-in, not a row insert/shift.

# Copy formatting from the row above (row 31) so number formats / alignment
# / wrap-text match the rest of the log (Sno/date/time/category/description).
$ws.Range("A31:G31").Copy() | Out-Null
$ws.Range("A32:G32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Populate the new entry (29 June 2022, 22:15 - 23:30, Code) ---
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 44741
$ws.Range("C32").Value = 0.92708333333333337
$ws.Range("D32").Value = 0.97916666666666663
$ws.Range("E32").Formula = "=D32-C32"
$ws.Range("F32").Value = "Code"
$ws.Range("G32").Value = "1. deeplab_v3_starter starter nb completed`n2. deeplab_v3_r50_baseline with r=[12, 24, 36] 10ep kaggle train`n3. deeplab_v3_r50_baseline with r=[6,12,18] 10ep kaggle train"

# Row holds a 3-line wrapped description -> matches the auto row height
# Excel computed for the other 3-line rows in this sheet (15pt x 3).
$ws.Rows("32").RowHeight = 45

# Totals formula (E33 = SUM(E2:E32)) already covers row 32, so no change
# is needed there; recalc runs automatically after the script completes and
# refreshes the displayed total to include the new entry.

# Reflect the last touched cell, same as the source edit.
$ws.Range("E33").Select() | Out-Null
